# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns populated, plus new hyperlinks on the
#    "Latest Target File" cells.
#  - Several columns are widened to fit the newly-populated long filenames.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetDisplay = "f3c8126a-19cc-4b2e-b2e7-9162bf5f5866.md"
$targetUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f2bd5dad278691b5bffc1a5f6d0d3786536dd04/e2e/f3c8126a-19cc-4b2e-b2e7-9162bf5f5866.md"

$zhHandbackFile = "f3c8126a-19cc-4b2e-b2e7-9162bf5f5866.42fe59d08295a0eb8a6f4d81701ff7622c75b5d0.zh-cn.xlf"
$deHandbackFile = "f3c8126a-19cc-4b2e-b2e7-9162bf5f5866.42fe59d08295a0eb8a6f4d81701ff7622c75b5d0.de-de.xlf"
$deHandbackDate = "2016-08-20 05:04:37"

# The host's ColumnWidth -> stored-width conversion snaps to 1/6-character
# pixel steps, so these "odd" inputs are chosen to round-trip to the wider
# widths (~29.98 and 40 "characters") used in the target layout.
$wideColWidth = 29.166666666666668
$maxColWidth  = 39.166666666666664

# --- Overview sheet: widen the zh-cn / de-de status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

function Update-LangSheet {
    param(
        [string]$SheetName,
        [string]$HandbackFile,
        [object]$HandbackDateValue
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) on rows 2 and 3.
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Latest Handback File (J) and Latest Handback DateTime (K).
    $ws.Range("J2").Value = $HandbackFile
    $ws.Range("J3").Value = $HandbackFile
    if ($HandbackDateValue -ne $null) {
        $ws.Range("K2").Value = $HandbackDateValue
        $ws.Range("K3").Value = $HandbackDateValue
    }

    # Widen columns C, I, J to fit the newly populated long values.
    $ws.Columns.Item(3).ColumnWidth = $wideColWidth
    $ws.Columns.Item(9).ColumnWidth = $maxColWidth
    $ws.Columns.Item(10).ColumnWidth = $maxColWidth

    # Latest Target File (I) - text + hyperlink. Leave the existing A2/A3
    # hyperlinks (handoff source files) untouched; only add the new ones.
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, "", "", $targetDisplay)
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, "", "", $targetDisplay)
}

Update-LangSheet "zh-cn" $zhHandbackFile $null

Update-LangSheet "de-de" $deHandbackFile $deHandbackDate

Write-Output "done"
